$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D10 from numeric 10 to a text value of three spaces
$ws.Range("D10").Value = "   "

# Add a new cell J17 with a text value of three spaces
$ws.Range("J17").Value = "   "

# Update the active selection to D10 (matches the saved selection in the diff)
$ws.Range("D10").Select()
